$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the typo in B1
$ws.Range("B1").Value = "Special Weekend"

# 2. Delete row 2 entirely - its date (45658) is dropped and all rows below
#    shift up by one (old row 3 becomes new row 2, etc.)
$ws.Rows(2).Delete()

# 3. The former C3 (now C2) held a stray "\n" placeholder; replace it with the
#    actual special-weekend timestamp text, and strip the leftover date
#    number format/style so it stores as plain text.
$ws.Range("C2").Value = "2025-04-26 00:00"
$ws.Range("C2").Style = "Normal"

# 4. Clear the stray "\n" placeholder text left in C3:C15 (after the shift)
$ws.Range("C3:C15").Value = ""

# 5. The last row (old row 17, now row 16) is fully blank - remove it so the
#    sheet's used range shrinks back down to row 15.
$ws.Rows(16).Delete()
